$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 1109; everything below shifts
# down by two rows (old row 1109 -> new row 1111, ..., old row 1172 -> new
# row 1174), matching the target dimension A1:T1174.
$ws.Rows("1109:1110").Insert()

# New row 1109: Murcott / Especial
$ws.Cells.Item(1109, 1).Value = 10
$ws.Cells.Item(1109, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1109, 3).Value = "La Araucanía"
$ws.Cells.Item(1109, 4).Value = 45267
$ws.Cells.Item(1109, 5).Value = 9
$ws.Cells.Item(1109, 6).Value = "Fruta"
$ws.Cells.Item(1109, 7).Value = 100102
$ws.Cells.Item(1109, 8).Value = "Cítricos"
$ws.Cells.Item(1109, 9).Value = 100102004
$ws.Cells.Item(1109, 10).Value = "Mandarina"
$ws.Cells.Item(1109, 11).Value = "Murcott"
$ws.Cells.Item(1109, 12).Value = "Especial"
$ws.Cells.Item(1109, 13).Value = 180
$ws.Cells.Item(1109, 14).Value = 18000
$ws.Cells.Item(1109, 15).Value = 18000
$ws.Cells.Item(1109, 16).Value = 18000
$ws.Cells.Item(1109, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(1109, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1109, 19).Value = 1000
$ws.Cells.Item(1109, 20).Value = 18

# New row 1110: Murcott / Primera
$ws.Cells.Item(1110, 1).Value = 10
$ws.Cells.Item(1110, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1110, 3).Value = "La Araucanía"
$ws.Cells.Item(1110, 4).Value = 45267
$ws.Cells.Item(1110, 5).Value = 9
$ws.Cells.Item(1110, 6).Value = "Fruta"
$ws.Cells.Item(1110, 7).Value = 100102
$ws.Cells.Item(1110, 8).Value = "Cítricos"
$ws.Cells.Item(1110, 9).Value = 100102004
$ws.Cells.Item(1110, 10).Value = "Mandarina"
$ws.Cells.Item(1110, 11).Value = "Murcott"
$ws.Cells.Item(1110, 12).Value = "Primera"
$ws.Cells.Item(1110, 13).Value = 155
$ws.Cells.Item(1110, 14).Value = 14000
$ws.Cells.Item(1110, 15).Value = 14000
$ws.Cells.Item(1110, 16).Value = 14000
$ws.Cells.Item(1110, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(1110, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1110, 19).Value = 778
$ws.Cells.Item(1110, 20).Value = 18
